$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 257-258 (everything from old row 257 onward shifts down by 2)
$ws.Range("A257:A258").EntireRow.Insert()

# New row 257: Camote, 1a (guarda), Región del Maule
$ws.Cells.Item(257,1).Value = 4
$ws.Cells.Item(257,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(257,3).Value = "Los Lagos"
$ws.Cells.Item(257,4).Value = 44694
$ws.Cells.Item(257,5).Value = 10
$ws.Cells.Item(257,6).Value = 100112045
$ws.Cells.Item(257,7).Value = "Zapallo"
$ws.Cells.Item(257,8).Value = "Camote"
$ws.Cells.Item(257,9).Value = "1a (guarda)"
$ws.Cells.Item(257,10).Value = 300
$ws.Cells.Item(257,11).Value = 600
$ws.Cells.Item(257,12).Value = 600
$ws.Cells.Item(257,13).Value = 600
$ws.Cells.Item(257,14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(257,15).Value = "Región del Maule"
$ws.Cells.Item(257,16).Value = 600
$ws.Cells.Item(257,17).Value = 1
$ws.Cells.Item(257,18).Value = "Hortaliza"

# New row 258: Paine, 1a (cosecha), Región de O'Higgins
$ws.Cells.Item(258,1).Value = 4
$ws.Cells.Item(258,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(258,3).Value = "Los Lagos"
$ws.Cells.Item(258,4).Value = 44694
$ws.Cells.Item(258,5).Value = 10
$ws.Cells.Item(258,6).Value = 100112045
$ws.Cells.Item(258,7).Value = "Zapallo"
$ws.Cells.Item(258,8).Value = "Paine"
$ws.Cells.Item(258,9).Value = "1a (cosecha)"
$ws.Cells.Item(258,10).Value = 800
$ws.Cells.Item(258,11).Value = 550
$ws.Cells.Item(258,12).Value = 550
$ws.Cells.Item(258,13).Value = 550
$ws.Cells.Item(258,14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(258,15).Value = "Región de O'Higgins"
$ws.Cells.Item(258,16).Value = 550
$ws.Cells.Item(258,17).Value = 1
$ws.Cells.Item(258,18).Value = "Hortaliza"

Write-Output "Inserted rows 257-258; new dimension rows:"
Write-Output ($ws.UsedRange.Rows.Count)
